$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns so numeric-looking
# strings (e.g. "521.24", "58.201.27") are stored as text, matching the
# original inlineStr cell type instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Column D (Price) / Column E (Volume 1h) value updates
$ws.Range("D2").Value = "58.201.27"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.593.82"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "521.24"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "144.32"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "2.615.92"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "0.102"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "3.053.31"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "58.184.77"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "20.54"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D19").Value = "341.32"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").Value = "4.37"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "10.33"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "6.40"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "66.13"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "0.404"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("D27").Value = "2.717.20"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "0.0₃0752"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "6.26"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "18.87"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "149.78"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "4.05"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "0.866"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "36.11"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "3.57"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "274.52"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "0.598"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "0.0958"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D50").Value = "19.19"
$ws.Range("E50").Value = "  +5.46%  "

# Row 17/18: ShibaInu and WrappedEther swapped position, with refreshed data
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.635.59"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("E18").Value = "  -0.45%  "

# Row 48/49: EnergySwap and Hedera swapped position, with refreshed data
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0526"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "18.84"
$ws.Range("E49").Value = "  -0.78%  "

# Row 51: RenderToken replaced by Maker
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.987.44"
$ws.Range("E51").Value = "  -2.10%  "

# Restore the cells to their original (default) formatting now that the
# text values are safely stored, so no extraneous number-format styling
# is left behind on the cells.
$ws.Range("D2:E51").ClearFormats()
